$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 35 for the new item "حلاوة حرير",
#    pushing "ديتول صغير" and everything below it down by one row.
$ws.Rows("35").Insert()

# 2. Copy cell formatting (number format / style) from the row that now sits
#    at 36 ("ديتول صغير") into the newly blank row 35, cell by cell, so the
#    new row matches the report's look (borders/fonts/number formats).
for ($c = 1; $c -le 17; $c++) {
    $src = $ws.Cells.Item(36, $c)
    $dst = $ws.Cells.Item(35, $c)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$ws.Rows("35").RowHeight = 24.75

# 3. Re-create the merged cell groups for the new row (A:B, C:G, H:K, L:M, N:O)
#    matching the pattern used by every other data row.
$ws.Range("A35:B35").Merge()
$ws.Range("C35:G35").Merge()
$ws.Range("H35:K35").Merge()
$ws.Range("L35:M35").Merge()
$ws.Range("N35:O35").Merge()

# 4. Make sure the text-like numeric columns stay stored as text (as every
#    other row in this report does), then fill in the new item's data.
$ws.Range("C35").NumberFormat = "@"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("L35").NumberFormat = "@"
$ws.Range("N35").NumberFormat = "@"
$ws.Range("P35").NumberFormat = "@"
$ws.Range("Q35").NumberFormat = "@"

$ws.Cells.Item(35, 1).Value = 29
$ws.Cells.Item(35, 3).Value = "حلاوة حرير"
$ws.Cells.Item(35, 8).Value = "9:0"
$ws.Cells.Item(35, 12).Value = "0"
$ws.Cells.Item(35, 14).Value = "20.00"
$ws.Cells.Item(35, 16).Value = "20.0000"
$ws.Cells.Item(35, 17).Value = "1:0"

# 5. Renumber the "م" (item #) column for every row that shifted down one
#    place, so the sequence stays 1..35 instead of jumping back to 29.
for ($r = 36; $r -le 41; $r++) {
    $cur = [double]$ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 1
}

# 6. Update the grand-total cell (now on row 42) to include the new item's
#    selling price (1768.81 + 20.00 = 1788.81).
$ws.Cells.Item(42, 16).Value = 1788.81

# 7. Update the generated-on timestamp in the footer (now on row 43).
$ws.Cells.Item(43, 1).Value = "Friday, 29 August, 2025 6:32 PM"

$ws.Cells.Item(1, 1).Select()
